# Updated IPS AIP hipo turnover
# Recomputed the "Internal Fill Rate" / turnover AOP (M:W, Jan..FY) monthly/quarterly/annual
# splits on several location tabs after the Professional & Manufacturing Voluntary Turnover
# AOP inputs (column E) changed. A few stray cells (that should stay blank like their
# neighbors) are cleared instead of holding a stale 0/placeholder value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Chino California")
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.0333333333333333
$ws.Range("P7").Value = 0.0333333333333333
$ws.Range("Q7").Value = 0.0333333333333333
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.0333333333333333
$ws.Range("T7").Value = 0.0333333333333333
$ws.Range("U7").Value = 0.0333333333333333
$ws.Range("V7").Value = 0.1
$ws.Range("W7").Value = 0.4

$ws = $wb.Worksheets.Item("El Paso Texas - EPC")
$ws.Range("M7").ClearContents()
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0.0112666666666667
$ws.Range("P10").Value = 0.0112666666666667
$ws.Range("Q10").Value = 0.0112666666666667
$ws.Range("R10").Value = 0.0338
$ws.Range("S10").Value = 0.0112666666666667
$ws.Range("T10").Value = 0.0112666666666667
$ws.Range("U10").Value = 0.0112666666666667
$ws.Range("V10").Value = 0.0338
$ws.Range("W10").Value = 0.1352

$ws = $wb.Worksheets.Item("Florence Kentucky")
$ws.Range("E7").Value = 0.067
$ws.Range("E8").Value = 0.067
$ws.Range("E9").Value = 0.067
$ws.Range("M9").Value = 0.0172
$ws.Range("N9").Value = 0.0343
$ws.Range("O9").Value = 0.0111666666666667
$ws.Range("P9").Value = 0.0111666666666667
$ws.Range("Q9").Value = 0.0111666666666667
$ws.Range("R9").Value = 0.0335
$ws.Range("S9").Value = 0.0111666666666667
$ws.Range("T9").Value = 0.0111666666666667
$ws.Range("U9").Value = 0.0111666666666667
$ws.Range("V9").Value = 0.0335
$ws.Range("W9").Value = 0.134

$ws = $wb.Worksheets.Item("Indianapolis Indiana")
$ws.Range("E2").Value = 0.0556
$ws.Range("E3").Value = 0.0556
$ws.Range("E4").Value = 0.0556
$ws.Range("M4").Value = 0.0556
$ws.Range("N4").Value = 0.0535
$ws.Range("O4").Value = 0.00926666666666667
$ws.Range("P4").Value = 0.00926666666666667
$ws.Range("Q4").Value = 0.00926666666666667
$ws.Range("R4").Value = 0.0278
$ws.Range("S4").Value = 0.00926666666666667
$ws.Range("T4").Value = 0.00926666666666667
$ws.Range("U4").Value = 0.00926666666666667
$ws.Range("V4").Value = 0.0278
$ws.Range("W4").Value = 0.1112
$ws.Range("M7").ClearContents()
$ws.Range("E8").Value = 0.186
$ws.Range("E9").Value = 0.186
$ws.Range("E10").Value = 0.186
$ws.Range("M10").Value = 0.0263
$ws.Range("N10").Value = 0.1241
$ws.Range("O10").Value = 0.031
$ws.Range("P10").Value = 0.031
$ws.Range("Q10").Value = 0.031
$ws.Range("R10").Value = 0.093
$ws.Range("S10").Value = 0.031
$ws.Range("T10").Value = 0.031
$ws.Range("U10").Value = 0.031
$ws.Range("V10").Value = 0.093
$ws.Range("W10").Value = 0.372

$ws = $wb.Worksheets.Item("Lavergne Tennessee")
$ws.Range("E2").Value = 0.069
$ws.Range("E3").Value = 0.069
$ws.Range("E4").Value = 0.069
$ws.Range("M4").Value = 0.0667
$ws.Range("N4").Value = 0.0667
$ws.Range("O4").Value = 0.0115
$ws.Range("P4").Value = 0.0115
$ws.Range("Q4").Value = 0.0115
$ws.Range("R4").Value = 0.0345
$ws.Range("S4").Value = 0.0115
$ws.Range("T4").Value = 0.0115
$ws.Range("U4").Value = 0.0115
$ws.Range("V4").Value = 0.0345
$ws.Range("W4").Value = 0.138
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0.5
$ws.Range("E8").Value = 0.0541
$ws.Range("E9").Value = 0.0541
$ws.Range("E10").Value = 0.0541
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0.00901666666666667
$ws.Range("P10").Value = 0.00901666666666667
$ws.Range("Q10").Value = 0.00901666666666667
$ws.Range("R10").Value = 0.02705
$ws.Range("S10").Value = 0.00901666666666667
$ws.Range("T10").Value = 0.00901666666666667
$ws.Range("U10").Value = 0.00901666666666667
$ws.Range("V10").Value = 0.02705
$ws.Range("W10").Value = 0.1082

$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("E2").Value = 0.3077
$ws.Range("E3").Value = 0.3077
$ws.Range("E4").Value = 0.3077
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0512833333333333
$ws.Range("P4").Value = 0.0512833333333333
$ws.Range("Q4").Value = 0.0512833333333333
$ws.Range("R4").Value = 0.15385
$ws.Range("S4").Value = 0.0512833333333333
$ws.Range("T4").Value = 0.0512833333333333
$ws.Range("U4").Value = 0.0512833333333333
$ws.Range("V4").Value = 0.15385
$ws.Range("W4").Value = 0.6154

$ws = $wb.Worksheets.Item("Pharr Texas")
$ws.Range("E2").Value = 0.1299
$ws.Range("E3").Value = 0.1299
$ws.Range("E4").Value = 0.1299
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.137
$ws.Range("O4").Value = 0.02165
$ws.Range("P4").Value = 0.02165
$ws.Range("Q4").Value = 0.02165
$ws.Range("R4").Value = 0.06495
$ws.Range("S4").Value = 0.02165
$ws.Range("T4").Value = 0.02165
$ws.Range("U4").Value = 0.02165
$ws.Range("V4").Value = 0.06495
$ws.Range("W4").Value = 0.2598
$ws.Range("M5").ClearContents()
$ws.Range("E6").Value = 0.0649
$ws.Range("E7").Value = 0.0649
$ws.Range("E8").Value = 0.0649
$ws.Range("M8").Value = 0.0312
$ws.Range("N8").Value = 0.0625
$ws.Range("O8").Value = 0.0108166666666667
$ws.Range("P8").Value = 0.0108166666666667
$ws.Range("Q8").Value = 0.0108166666666667
$ws.Range("R8").Value = 0.03245
$ws.Range("S8").Value = 0.0108166666666667
$ws.Range("T8").Value = 0.0108166666666667
$ws.Range("U8").Value = 0.0108166666666667
$ws.Range("V8").Value = 0.03245
$ws.Range("W8").Value = 0.1298

Write-Host "done"